$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 661 entirely ("「開⇔閉」" entry); all following rows shift up by one.
$ws.Rows.Item(661).Delete()
